$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# coinranking.com snapshot refresh (GitHub Actions cron).
# Price/Volume columns are stored as *text* (mixed "1.234.56"-style
# thousands-dot formatting, scientific-looking decimals, "  +1.23%  " with
# padding, ...). Plain `.Value = "..."` lets Excel's smart-entry parse
# number-shaped strings (e.g. "244.74", "1.000") into real numbers, so for
# those cells we force Text format, assign, then restore the default style
# (so no stray number format sticks to the cell afterwards).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "29.353.72"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "1.846.41"
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "244.74"
Set-TextValue $ws.Range("D6") "0.6908"
$ws.Range("E6").Value = "  -0.07%  "
Set-TextValue $ws.Range("D8") "0.3052"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.27%  "
Set-TextValue $ws.Range("D10") "23.39"
$ws.Range("E10").Value = "  +0.23%  "
Set-TextValue $ws.Range("D11") "0.07718"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D12") "5.131"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.846.35"
$ws.Range("E13").Value = "  +0.66%  "
Set-TextValue $ws.Range("D14") "0.6896"
$ws.Range("E14").Value = "  +1.32%  "
Set-TextValue $ws.Range("D15") "90.07"
$ws.Range("E15").Value = "  -3.44%  "
Set-TextValue $ws.Range("D16") "6.304"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").Value = "29.350.51"
$ws.Range("E17").Value = "  +1.60%  "
Set-TextValue $ws.Range("D18") "0.000008237"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "2.093.74"
$ws.Range("E19").Value = "  +0.56%  "
Set-TextValue $ws.Range("D20") "235.96"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("E22").Value = "  +0.04%  "
Set-TextValue $ws.Range("D23") "7.639"
$ws.Range("E23").Value = "  +2.51%  "
Set-TextValue $ws.Range("D24") "1.000"
$ws.Range("E24").Value = "  +0.03%  "
Set-TextValue $ws.Range("D25") "0.1471"
$ws.Range("E25").Value = "  -0.40%  "
Set-TextValue $ws.Range("D26") "8.930"
$ws.Range("E26").Value = "  +1.52%  "
Set-TextValue $ws.Range("D27") "160.13"
$ws.Range("E27").Value = "  +0.37%  "
Set-TextValue $ws.Range("D28") "18.13"
$ws.Range("E28").Value = "  -0.80%  "
Set-TextValue $ws.Range("D29") "1.526"
$ws.Range("E29").Value = "  -0.92%  "
Set-TextValue $ws.Range("D30") "4.247"
$ws.Range("E30").Value = "  +0.62%  "
Set-TextValue $ws.Range("D31") "4.134"
$ws.Range("E31").Value = "  -0.62%  "
Set-TextValue $ws.Range("D32") "1.201"
$ws.Range("E32").Value = "  +1.16%  "
Set-TextValue $ws.Range("D33") "0.05217"
$ws.Range("E33").Value = "  +2.37%  "
Set-TextValue $ws.Range("D34") "0.7714"
$ws.Range("E34").Value = "  -0.43%  "
Set-TextValue $ws.Range("D35") "1.873"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("E36").Value = "  +0.14%  "
Set-TextValue $ws.Range("D37") "2.678"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Value = "1.304.82"
$ws.Range("E38").Value = "  +5.79%  "
Set-TextValue $ws.Range("D39") "0.01860"
$ws.Range("E39").Value = "  +0.53%  "
Set-TextValue $ws.Range("D40") "2.706"
$ws.Range("E40").Value = "  +0.32%  "
Set-TextValue $ws.Range("D41") "0.9440"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  -2.23%  "
Set-TextValue $ws.Range("D43") "5.752"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("E44").Value = "  +0.01%  "
Set-TextValue $ws.Range("D45") "9.713"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "1.994.40"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("E47").Value = "  +1.17%  "
Set-TextValue $ws.Range("D48") "1.777"
$ws.Range("E48").Value = "  +1.91%  "
Set-TextValue $ws.Range("D49") "0.00000000120"
$ws.Range("E49").Value = "  +2.05%  "
Set-TextValue $ws.Range("D50") "63.30"
$ws.Range("E50").Value = "  -1.13%  "
Set-TextValue $ws.Range("D51") "0.05936"
$ws.Range("E51").Value = "  +0.79%  "
